# Add a new data row (row 3) to the "Artfynd" sheet, mirroring the
# structure of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some columns in this export hold date/time-looking strings that must stay
# literal text (the source file stores them as inline strings, not real
# Excel dates). Force a text number format on those cells before writing so
# Excel's autodetection doesn't convert them into date serial numbers.
$ws.Range("Y3:AB3").NumberFormat = "@"

$ws.Range("A3").Value = 111833710
$ws.Range("B3").Value = 88966
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 5754
$ws.Range("F3").Value = "Gultoppig fingersvamp"
$ws.Range("G3").Value = "Ramaria testaceoflava"
$ws.Range("H3").Value = "(Bres.) Corner"
$ws.Range("I3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("P3").Value = "Brudslöjan (Brudslöjan), Hjd"
$ws.Range("Q3").Value = 447323.2079976397
$ws.Range("R3").Value = 6929278.806905948
$ws.Range("S3").Value = 25
$ws.Range("T3").Value = "Jämtland"
$ws.Range("U3").Value = "Härjedalen"
$ws.Range("V3").Value = "Härjedalen"
$ws.Range("W3").Value = "Vemdalen"
$ws.Range("Y3").Value = "2023-09-01"
$ws.Range("Z3").Value = "10:25"
$ws.Range("AA3").Value = "2023-09-01"
$ws.Range("AB3").Value = "10:25"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AT3").Value = ""
$ws.Range("AW3").Value = "Tommy Carlström"
$ws.Range("AX3").Value = "Tommy Carlström, Ingela Carlström, Thomas Samuelsson, Ylva Rinaldo"
$ws.Range("AY3").Value = ""
